# tests(US-Tool): add new test
# Adds a new "central_h2_heat_link" row to the 6_links sheet and a new
# "central_h2_heat_bus" row to the 1_buses sheet (inserted in the middle
# of the existing bus list), matching the standard_parameters.xlsx test
# fixture update.

$wb = $excel.ActiveWorkbook

# --- 1) 6_links: append the new link row (row 20) -------------------------
# NOTE: this has to happen before the 1_buses edit below so that the new
# shared-string "central_h2_heat_link" gets interned before
# "central_h2_heat_bus" (matches shared string order 330/331).
$ws9 = $wb.Worksheets.Item("6_links")

$ws9.Range("A20").Value = "central_h2_heat_link"
$ws9.Range("B20").Value = 1
$ws9.Range("C20").Value = "directed"
$ws9.Range("D20").Value = 1
$ws9.Range("E20").Value = 0
$ws9.Range("F20").Value = 0
$ws9.Range("G20").Value = 9999
$ws9.Range("H20").Value = 0
$ws9.Range("I20").Value = 0
$ws9.Range("J20").Value = 0.00001
$ws9.Range("K20").Value = 0.00001
$ws9.Range("L20").Value = 0
$ws9.Range("M20").Value = 0
$ws9.Range("N20").Value = 0

# conditional formatting for the new row, split the same way the old
# "last row" (row 19) rule was split: label columns vs. numeric columns.
$fc9a = $ws9.Range("A20:C20").FormatConditions.Add(1, 3, "0")
$fc9a.Interior.Color = 8421504
$fc9b = $ws9.Range("D20:N20").FormatConditions.Add(1, 3, "0")
$fc9b.Interior.Color = 8421504

# --- 2) 1_buses: insert the new bus row (row 23) ---------------------------
$ws2 = $wb.Worksheets.Item("1_buses")

$ws2.Rows("23:23").Insert()
$ws2.Range("A23").Value = "central_h2_heat_bus"
$ws2.Range("B23").Value = 1
$ws2.Range("C23").Value = 1
$ws2.Range("D23").Value = 0
$ws2.Range("E23").Value = 0
$ws2.Range("F23").Value = 0
$ws2.Range("G23").Value = 0
$ws2.Range("H23").Value = 0

# the newly inserted row keeps no conditional formatting (same as the
# original file), but the row insert leaves the old single rule covering
# B9:H34 untouched/un-shifted, so rebuild it excluding the new row 23.
$oldRule = $ws2.Range("B9:H34")
if ($oldRule.FormatConditions.Count -gt 0) {
    $oldRule.FormatConditions.Item(1).Delete()
}

$fc2a = $ws2.Range("B3:H5").FormatConditions.Add(1, 3, "0")
$fc2a.Interior.Color = 8421504
$fc2b = $ws2.Range("B9:H22").FormatConditions.Add(1, 3, "0")
$fc2b.Interior.Color = 8421504
$fc2c = $ws2.Range("F6:F6").FormatConditions.Add(1, 3, "0")
$fc2c.Interior.Color = 8421504
$fc2d = $ws2.Range("B24:H35").FormatConditions.Add(1, 3, "0")
$fc2d.Interior.Color = 8421504

# --- 3) leftover selection/active-cell bookkeeping -------------------------
$ws7 = $wb.Worksheets.Item("4_transformers")
$ws7.Range("B29").Select()

$ws9.Range("C29").Select()

$ws2.Activate()
$ws2.Range("L29").Select()
